$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.209624
$ws.Range("H2").Value = 153.628872
$ws.Range("I2").Value = 0.9420574924009609
$ws.Range("J2").Value = 0.9606108937376658
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 19.5719925
$ws.Range("N2").Value = 39.143985
$ws.Range("O2").Value = 0.03094210933382397
$ws.Range("P2").Value = 0.02187190777676379
$ws.Range("Q2").Value = 1002.27437685582
$ws.Range("R2").Value = 6013.64626113492
$ws.Range("S2").Value = 0.02914924592861858
$ws.Range("T2").Value = 0.02101039287718487
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.209624
$ws.Range("H3").Value = 153.628872
$ws.Range("I3").Value = 0.9420574924009609
$ws.Range("J3").Value = 0.9606108937376658
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 115.495743
$ws.Range("N3").Value = 346.487229
$ws.Range("O3").Value = 0.1825916246134488
$ws.Range("P3").Value = 0.1936015640337701
$ws.Range("Q3").Value = 5914.493572630632
$ws.Range("R3").Value = 53230.44215367569
$ws.Range("S3").Value = 0.1720118080167631
$ws.Range("T3").Value = 0.1859757714554899
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.209624
$ws.Range("H4").Value = 153.628872
$ws.Range("I4").Value = 0.9420574924009609
$ws.Range("J4").Value = 0.9606108937376658
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 239.8982746666667
$ws.Range("N4").Value = 719.694824
$ws.Range("O4").Value = 0.3792643310961689
$ws.Range("P4").Value = 0.4021332732970914
$ws.Range("Q4").Value = 12285.10044392873
$ws.Range("R4").Value = 110565.9039953585
$ws.Range("S4").Value = 0.3572888047095846
$ws.Range("T4").Value = 0.386293603063572
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.209624
$ws.Range("H5").Value = 153.628872
$ws.Range("I5").Value = 0.9420574924009609
$ws.Range("J5").Value = 0.9606108937376658
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 151.102183
$ws.Range("N5").Value = 453.306549
$ws.Range("O5").Value = 0.2388832034840335
$ws.Range("P5").Value = 0.2532874216646837
$ws.Range("Q5").Value = 7737.885977009191
$ws.Range("R5").Value = 69640.97379308273
$ws.Range("S5").Value = 0.2250417116508771
$ws.Range("T5").Value = 0.2433106564978209
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 51.209624
$ws.Range("H6").Value = 153.628872
$ws.Range("I6").Value = 0.9420574924009609
$ws.Range("J6").Value = 0.9606108937376658
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.12446233333333
$ws.Range("N6").Value = 54.373387
$ws.Range("O6").Value = 0.02865365369084289
$ws.Range("P6").Value = 0.03038141635232813
$ws.Range("Q6").Value = 928.1469012921626
$ws.Range("R6").Value = 8353.322111629464
$ws.Range("S6").Value = 0.02699338914412099
$ws.Range("T6").Value = 0.02918471951522606
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 51.209624
$ws.Range("H7").Value = 153.628872
$ws.Range("I7").Value = 0.9420574924009609
$ws.Range("J7").Value = 0.9606108937376658
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 88.3431645
$ws.Range("N7").Value = 176.686329
$ws.Range("O7").Value = 0.139665077781682
$ws.Range("P7").Value = 0.09872441687536272
$ws.Range("Q7").Value = 4524.020237015148
$ws.Range("R7").Value = 27144.12142209089
$ws.Range("S7").Value = 0.1315725329509965
$ws.Range("T7").Value = 0.09483575032837209
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 3.1497165
$ws.Range("H8").Value = 6.299433000000001
$ws.Range("I8").Value = 0.05794250759903903
$ws.Range("J8").Value = 0.03938910626233424
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 19.5719925
$ws.Range("N8").Value = 39.143985
$ws.Range("O8").Value = 0.03094210933382397
$ws.Range("P8").Value = 0.02187190777676379
$ws.Range("Q8").Value = 61.64622771512626
$ws.Range("R8").Value = 246.584910860505
$ws.Range("S8").Value = 0.001792863405205392
$ws.Range("T8").Value = 0.0008615148995789236
$ws.Range("E9").Value = 2
$ws.Range("G9").Value = 3.1497165
$ws.Range("H9").Value = 6.299433000000001
$ws.Range("I9").Value = 0.05794250759903903
$ws.Range("J9").Value = 0.03938910626233424
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 115.495743
$ws.Range("N9").Value = 346.487229
$ws.Range("O9").Value = 0.1825916246134488
$ws.Range("P9").Value = 0.1936015640337701
$ws.Range("Q9").Value = 363.7788474068595
$ws.Range("R9").Value = 2182.673084441157
$ws.Range("S9").Value = 0.01057981659668564
$ws.Range("T9").Value = 0.007625792578280278
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 3.1497165
$ws.Range("H10").Value = 6.299433000000001
$ws.Range("I10").Value = 0.05794250759903903
$ws.Range("J10").Value = 0.03938910626233424
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 239.8982746666667
$ws.Range("N10").Value = 719.694824
$ws.Range("O10").Value = 0.3792643310961689
$ws.Range("P10").Value = 0.4021332732970914
$ws.Range("Q10").Value = 755.6115540391321
$ws.Range("R10").Value = 4533.669324234793
$ws.Range("S10").Value = 0.02197552638658422
$ws.Range("T10").Value = 0.01583967023351943
$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 3.1497165
$ws.Range("H11").Value = 6.299433000000001
$ws.Range("I11").Value = 0.05794250759903903
$ws.Range("J11").Value = 0.03938910626233424
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 151.102183
$ws.Range("N11").Value = 453.306549
$ws.Range("O11").Value = 0.2388832034840335
$ws.Range("P11").Value = 0.2532874216646837
$ws.Range("Q11").Value = 475.9290389811195
$ws.Range("R11").Value = 2855.574233886717
$ws.Range("S11").Value = 0.0138414918331564
$ws.Range("T11").Value = 0.009976765166862888
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 3.1497165
$ws.Range("H12").Value = 6.299433000000001
$ws.Range("I12").Value = 0.05794250759903903
$ws.Range("J12").Value = 0.03938910626233424
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.12446233333333
$ws.Range("N12").Value = 54.373387
$ws.Range("O12").Value = 0.02865365369084289
$ws.Range("P12").Value = 0.03038141635232813
$ws.Range("Q12").Value = 57.0869180649285
$ws.Range("R12").Value = 342.521508389571
$ws.Range("S12").Value = 0.001660264546721897
$ws.Range("T12").Value = 0.001196696837102072
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 3.1497165
$ws.Range("H13").Value = 6.299433000000001
$ws.Range("I13").Value = 0.05794250759903903
$ws.Range("J13").Value = 0.03938910626233424
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 88.3431645
$ws.Range("N13").Value = 176.686329
$ws.Range("O13").Value = 0.139665077781682
$ws.Range("P13").Value = 0.09872441687536272
$ws.Range("Q13").Value = 278.2559228878643
$ws.Range("R13").Value = 1113.023691551457
$ws.Range("S13").Value = 0.008092544830685484
$ws.Range("T13").Value = 0.003888666546990646
